$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample number text from "E7420" to "E7420L" for the whole
# roboticS2Prep column (G2:G41) which all shared the same string value.
$ws.Range("G2:G41").Value = "E7420L"

# These cells previously held a "=FALSE()" formula evaluating to FALSE;
# replace them with a plain boolean constant (no formula) of FALSE.
$ws.Range("H2:H41").Value = $false

# Update the sheet's active selection to just G2 (was G2:G41).
[void]$ws.Range("G2").Select()
